$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "75.608.40"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.870.44"
$ws.Range("E3").Value = "  +6.51%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "194.76"
$ws.Range("E5").Value = "  +3.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "596.63"
$ws.Range("E6").Value = "  +1.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.551"
$ws.Range("E8").Value = "  +1.79%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.191"
$ws.Range("E9").Value = "  -3.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "2.876.51"
$ws.Range("E10").Value = "  +6.82%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.397"
$ws.Range("E11").Value = "  +10.34%  "
$ws.Range("E12").Value = "  -1.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.89"
$ws.Range("E13").Value = "  +2.96%  "
$ws.Range("E14").Value = "  +7.91%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "75.679.47"
$ws.Range("E15").Value = "  +0.19%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000188"
$ws.Range("E16").Value = "  -0.78%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.20"
$ws.Range("E17").Value = "  +1.74%  "
$ws.Range("E18").Value = "  +7.74%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.93"
$ws.Range("E19").Value = "  -4.80%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.49"
$ws.Range("E20").Value = "  +3.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "379.98"
$ws.Range("E21").Value = "  +1.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.28"
$ws.Range("E22").Value = "  -0.72%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.13"
$ws.Range("E23").Value = "  +1.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.78"
$ws.Range("E24").Value = "  +1.32%  "
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.20"
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.69"
$ws.Range("E28").Value = "  +1.98%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000106"
$ws.Range("E29").Value = "  +10.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -0.26%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.39"
$ws.Range("E31").Value = "  -1.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "506.24"
$ws.Range("E32").Value = "  -3.35%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.73"
$ws.Range("E33").Value = "  -1.65%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.81"
$ws.Range("E34").Value = "  +1.92%  "
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "165.74"
$ws.Range("E36").Value = "  +1.91%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.05"
$ws.Range("E37").Value = "  +3.51%  "
$ws.Range("B38").Value = "WhiteBITCoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.68"
$ws.Range("E38").Value = "  +1.47%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.114"
$ws.Range("E39").Value = "  -5.39%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "182.57"
$ws.Range("E40").Value = "  +6.02%  "
$ws.Range("E41").Value = "  -0.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.343"
$ws.Range("E42").Value = "  +2.93%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.98"
$ws.Range("E43").Value = "  -1.66%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.67"
$ws.Range("E44").Value = "  -2.51%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0908"
$ws.Range("E45").Value = "  +6.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.21"
$ws.Range("E46").Value = "  +0.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "40.19"
$ws.Range("E47").Value = "  +1.99%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.33"
$ws.Range("E48").Value = "  -3.14%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.578"
$ws.Range("E49").Value = "  +6.40%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.75"
$ws.Range("E50").Value = "  +1.94%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.661"
$ws.Range("E51").Value = "  +11.20%  "
